$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of results (puntajes) to append below the existing header row.
# Columns: A=PUNT_INGLES B=PUNT_MATEMATICAS C=PUNT_SOCIALES_CIUDADANAS
#          D=PUNT_C_NATURALES E=PUNT_LECTURA_CRITICA F=puntaje_global G=anno
$data = @(
    @(2, 0.580118231871416, 0.662182702606799, 0,                 0,                 0,                 0.197435872283986, "20122"),
    @(3, 0.580055201794936, 0.626663972406649, 0,                 0,                 0,                 0.189234393770375, "20132"),
    @(4, 0.665790928107443, 0.673954640202497, 0.658527350237695, 0.652508764251948, 0.684904529089117, 0.667344367650093, "20142"),
    @(5, 0.679211515958017, 0.752438171964393, 0.704808118809078, 0.680656520337941, 0.703394916032865, 0.707931130568527, "20152"),
    @(6, 0.733530459943224, 0.765028054660648, 0.720420186842568, 0.736529564848269, 0.736092012106156, 0.73905691656278,  "20162"),
    @(7, 0.695613509580673, 0.719358733300511, 0.706520143068675, 0.709106184409968, 0.725890402607448, 0.713710761518499, "20172"),
    @(8, 0.672230939908778, 0.695073585543965, 0.655085660754569, 0.673462481001246, 0.727203793947279, 0.686515961819228, "20194"),
    @(9, 0.696695863693876, 0.680398749386295, 0.664501687680721, 0.661505967723753, 0.716338672460317, 0.681917776572856, "20224")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]

    # "anno" is a text value (e.g. "20122") even though it looks numeric, so
    # force it to be stored as text (leading apostrophe), then strip the
    # resulting "quote prefix" cell format back to Normal so no stray
    # number-format is left behind on the cell.
    $ws.Range("G$r").Value = "'" + $row[7]
    $ws.Range("G$r").Style = "Normal"
}
